# Apply the "Approved" column rework + reject-page related cleanup
# to the Items worksheet.
#
# Summary of the change (per the commit's diff):
#  - Rows 9, 10 and 11 (items 8, 9, 10) are removed from the sheet.
#  - The J/K header columns are swapped: J1 becomes "Approved" and
#    K1 becomes "ReservedBy" (previously J1="ReservedBy", K1="Approved").
#  - For every remaining data row (2-8) the J column now holds a
#    boolean TRUE ("Approved"), and whatever e-mail address used to
#    live in J (the old "ReservedBy" value) is moved over to K.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Items")

# --- Remove the last three rows (items 8, 9 and 10) -----------------
$ws.Rows.Item(9).Resize(3).EntireRow.Delete() | Out-Null

# --- Capture the old J-column (ReservedBy) values before overwriting -
$reservedBy = @{}
for ($r = 2; $r -le 8; $r++) {
    $val = $ws.Cells.Item($r, 10).Value()
    if ($null -ne $val -and $val -ne "") {
        $reservedBy[$r] = $val
    }
}

# --- Swap the header labels in J1 / K1 -------------------------------
$ws.Cells.Item(1, 10).Value = "Approved"
$ws.Cells.Item(1, 11).Value = "ReservedBy"

# --- Update each data row: J = Approved (TRUE), K = ReservedBy email -
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 10).Value = $true
    if ($reservedBy.ContainsKey($r)) {
        $ws.Cells.Item($r, 11).Value = $reservedBy[$r]
    }
}
